$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "family" column is column A; the relative-abundance sample columns
# (B..Y) each sum to 100 across all family rows. We are filtering out the
# "Partitiviridae" family row entirely, which requires:
#   1. Removing its row (shared string + row entry) so everything below
#      shifts up.
#   2. Rescaling column G (the only column with a nonzero contribution
#      from the removed row) so it again sums to 100.

# Locate the row containing the "Partitiviridae" family name in column A.
$found = $ws.Columns.Item(1).Find("Partitiviridae")
$targetRow = $found.Row

# Column G is the 7th column; remember the value being removed so we can
# rescale the remaining rows proportionally.
$removedG = $ws.Cells.Item($targetRow, 7).Value2

# Determine the full extent of the data (last used row/column).
$lastRow = $ws.UsedRange.Rows.Count
$lastCol = $ws.UsedRange.Columns.Count

# Remove the entire row; Excel will shift all subsequent rows up and drop
# the now-unused shared string automatically.
$ws.Rows.Item($targetRow).Delete()

$newLastRow = $lastRow - 1

# Rescale every remaining data row's column G value so the column totals
# 100 again, matching the other (already-100-summing) columns.
$scale = 100.0 / (100.0 - $removedG)

for ($r = 2; $r -le $newLastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $current = $cell.Value2
    $cell.Value = $current * $scale
}
